$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, formatted the same way as the other headers (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill H2:H50 with 0 (no-format change needed, they match the plain default style)
$ws.Range("H2:H50").Value = 0

# Row 13 is the only save in this data set
$ws.Range("H13").Value = 1
